$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It is inserted as
# row 49 (date 2021-10-19), pushing the previously-existing rows 49-74
# down to become rows 50-75.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record's data.
$ws.Cells.Item(49, 1).Value = 4
$ws.Cells.Item(49, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(49, 3).Value = 'Los Lagos'
$ws.Cells.Item(49, 4).Value = 44488
$ws.Cells.Item(49, 5).Value = 10
$ws.Cells.Item(49, 6).Value = 100112022
$ws.Cells.Item(49, 7).Value = 'Arveja Verde'
$ws.Cells.Item(49, 8).Value = 'Sin especificar'
$ws.Cells.Item(49, 9).Value = 'Primera'
$ws.Cells.Item(49, 10).Value = 80
$ws.Cells.Item(49, 11).Value = 25000
$ws.Cells.Item(49, 12).Value = 25000
$ws.Cells.Item(49, 13).Value = 25000
$ws.Cells.Item(49, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(49, 15).Value = 'Región del Maule'
$ws.Cells.Item(49, 16).Value = 1000
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = 'Hortaliza'
